$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A14 loses its underline formatting (style index 2 -> default style)
$ws.Range("A14").Font.Underline = 0

# New menu items (rows 70-78)
$ws.Range("A70").Value = "Iogurte Grego"
$ws.Range("B70").Value = "Lanche Manhã"
$ws.Range("D70").Value = 129

$ws.Range("A71").Value = "Iogurte Grego"
$ws.Range("B71").Value = "Lanche Tarde"
$ws.Range("D71").Value = 129

$ws.Range("A72").Value = "Iogurte Grego"
$ws.Range("B72").Value = "Café da Manhã"
$ws.Range("D72").Value = 129

$ws.Range("A73").Value = "Crackers Pingo Doce"
$ws.Range("B73").Value = "Café da Manhã"
$ws.Range("D73").Value = 437

$ws.Range("A74").Value = "Crackers Pingo Doce"
$ws.Range("B74").Value = "Lanche Manhã"
$ws.Range("D74").Value = 437

$ws.Range("A75").Value = "Crackers Pingo Doce"
$ws.Range("B75").Value = "Lanche Tarde"
$ws.Range("D75").Value = 437

$ws.Range("A76").Value = "Crackers Mercadona"
$ws.Range("B76").Value = "Café da Manhã"
$ws.Range("D76").Value = 443

$ws.Range("A77").Value = "Crackers Mercadona"
$ws.Range("B77").Value = "Lanche Manhã"
$ws.Range("D77").Value = 443

$ws.Range("A78").Value = "Crackers Mercadona"
$ws.Range("B78").Value = "Lanche Tarde"
$ws.Range("D78").Value = 443

# Restore selection / scroll state to match the saved view
$ws.Range("A14").Select()
